$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 9 with new data (previously blank row)
$ws.Range("A9").Value = "Week2.5"
$ws.Range("B9").Value = "Other players knows which card is suppressed"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = "Need to create a new random code for each new game for each card, different from the card ID"

# Update the selected cell
$ws.Range("F14").Select()
